# Updates cryptos list figures (prices / 1h volume change %) and fixes the
# ordering of a couple of rows whose coin name/link/price/volume had been
# mixed up (rows 37/38, 45/46) plus swaps Monero -> Mantle in row 51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Some of the new "Price" values look like plain numbers (e.g. "0.625").
# The source file stores every Price/Volume cell as text, so force the
# cells we are about to overwrite with number-looking text to stay text
# by pre-formatting them with a text number format. Cells whose new value
# still contains extra punctuation (e.g. "69.901.83") already round-trip
# as text automatically and do not need this treatment, but it is
# harmless to include them too.
# NOTE: applying NumberFormat to a multi-area Union range only affects
# the first area in this engine, so the cells are set individually.
$cellsNeedingTextFormat = @("D2","D3","D5","D6","D7","D10","D11","D13","D14","D15","D16","D17","D18","D19","D21","D22","D24","D25","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D50","D51")
foreach ($addr in $cellsNeedingTextFormat) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.901.83"

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.519.74"
$ws.Range("E3").Value = "  -0.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
$ws.Range("D5").Value = "601.65"
$ws.Range("E5").Value = "  -1.63%  "

# Row 6 - Solana
$ws.Range("D6").Value = "196.03"
$ws.Range("E6").Value = "  +6.04%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  +0.90%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -2.04%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "0.654"
$ws.Range("E10").Value = "  +1.45%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "54.14"
$ws.Range("E11").Value = "  +1.27%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -2.41%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "9.55"
$ws.Range("E13").Value = "  +1.49%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.075.98"
$ws.Range("E14").Value = "  -0.62%  "

# Row 15 - BitcoinCash
$ws.Range("D15").Value = "603.84"
$ws.Range("E15").Value = "  -0.96%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "70.068.91"
$ws.Range("E16").Value = "  +0.22%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "19.12"
$ws.Range("E17").Value = "  +1.61%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "12.60"
$ws.Range("E18").Value = "  -0.21%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.514.55"
$ws.Range("E19").Value = "  -1.00%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.72%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "0.994"
$ws.Range("E21").Value = "  +0.37%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "18.22"
$ws.Range("E22").Value = "  +3.99%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  +6.47%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "104.00"
$ws.Range("E24").Value = "  +3.80%  "

# Row 25
$ws.Range("D25").Value = "4.60"
$ws.Range("E25").Value = "  -2.58%  "

# Row 26
$ws.Range("E26").Value = "  +2.85%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  -0.15%  "

# Row 28 - Filecoin
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "33.57"
$ws.Range("E29").Value = "  +3.50%  "

# Row 30 - dogwifhat
$ws.Range("D30").Value = "4.51"
$ws.Range("E30").Value = "  +25.56%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  +1.70%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "12.74"
$ws.Range("E32").Value = "  +4.32%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +1.93%  "

# Row 34 - OKB
$ws.Range("D34").Value = "63.19"
$ws.Range("E34").Value = "  -0.32%  "

# Row 35 - Maker
$ws.Range("D35").Value = "3.756.04"
$ws.Range("E35").Value = "  +6.40%  "

# Row 36 - PEPE
$ws.Range("D36").Value = "0.0₃0819"
$ws.Range("E36").Value = "  +5.23%  "

# Row 37 - was Dai, now Fetch.AI
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "3.09"
$ws.Range("E37").Value = "  -4.55%  "

# Row 38 - was Fetch.AI, now Dai
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.04%  "

# Row 39 - TheGraph
$ws.Range("D39").Value = "0.394"
$ws.Range("E39").Value = "  -1.09%  "

# Row 40 - Stacks
$ws.Range("D40").Value = "3.61"
$ws.Range("E40").Value = "  +1.29%  "

# Row 41 - InjectiveProtocol
$ws.Range("D41").Value = "36.85"
$ws.Range("E41").Value = "  -0.58%  "

# Row 42 - Bittensor
$ws.Range("D42").Value = "491.48"
$ws.Range("E42").Value = "  -7.97%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -0.18%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0457"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45 - was ApeXProtocol, now Stellar
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.141"
$ws.Range("E45").Value = "  -0.92%  "

# Row 46 - was Stellar, now ApeXProtocol
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.33"
$ws.Range("E46").Value = "  -0.72%  "

# Row 47 - ThetaToken
$ws.Range("E47").Value = "  -2.81%  "

# Row 48 - FirstDigitalUSD
$ws.Range("E48").Value = "  +0.37%  "

# Row 49 - THORChain
$ws.Range("E49").Value = "  -5.53%  "

# Row 50 - FLOKI
$ws.Range("D50").Value = "0.000243"
$ws.Range("E50").Value = "  +1.00%  "

# Row 51 - was Monero, now Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "1.31"
$ws.Range("E51").Value = "  +11.47%  "
